$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1. Latest HO Xliff Generate Date / Correspond Handoff Datetime (de-de) - shared string "2016-11-09 07:04:18" -> "2016-11-09 07:06:20"
$wsOverview.Range("G2").Value = "2016-11-09 07:06:20"
$wsOverview.Range("G3").Value = "2016-11-09 07:06:20"
$wsDeDe.Range("H2").Value = "2016-11-09 07:06:20"
$wsDeDe.Range("H3").Value = "2016-11-09 07:06:20"

# 2. Priority column "ht" -> "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# 3. zh-cn Correspond Handoff Datetime "2016-11-09 07:03:56" -> "2016-11-09 07:06:05"
$wsZhCn.Range("H2").Value = "2016-11-09 07:06:05"
$wsZhCn.Range("H3").Value = "2016-11-09 07:06:05"

# 4. zh-cn Correspond Handback DateTime "2016-11-09 07:05:03" -> "2016-11-09 07:06:59"
$wsZhCn.Range("K2").Value = "2016-11-09 07:06:59"
$wsZhCn.Range("K3").Value = "2016-11-09 07:06:59"

# 5. de-de Correspond Handback DateTime "2016-11-09 07:05:23" -> "2016-11-09 07:07:17"
$wsDeDe.Range("K2").Value = "2016-11-09 07:07:17"
$wsDeDe.Range("K3").Value = "2016-11-09 07:07:17"
